$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.031.78'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.485.80'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.74'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.89'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.618'
$ws.Range('E7').Value = '  +3.83%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '3.484.55'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').Value = '4.087.80'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.20'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '68.015.28'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000178'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '3.489.35'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.16'
$ws.Range('E20').Value = '  -4.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '396.01'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.96'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.29'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').Value = '  -1.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.47'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.41'
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.63'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  -4.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.00'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.893'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('E40').Value = '  +5.88%  '
$ws.Range('E41').Value = '  -4.74%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.80'
$ws.Range('E42').Value = '  -5.28%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.69'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.24'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0720'
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.39'
$ws.Range('E46').Value = '  -4.79%  '
$ws.Range('D47').Value = '2.755.18'
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.52'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0299'
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '330.56'
$ws.Range('E50').Value = '  -4.49%  '
$ws.Range('E51').Value = '  -3.58%  '
